# Apply the price / coin-ranking updates scraped on 2022-12-18.
# Column D holds text-formatted numbers (e.g. "6.470"), so force a
# text number format before writing to avoid Excel silently
# re-typing the value as a Number and dropping trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '246.41'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.537'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05634'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.470'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8064'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.058'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1432'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07301'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03226'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02927'
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09261'
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001666'
$ws.Range("E14").Value = '13BitForexTokenBF'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.216'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04715'
$ws.Range("E16").Value = '15CoinExTokenCET'
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0005828'
$ws.Range("E17").Value = '16OneONE'
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006274'
$ws.Range("E18").Value = '17TigerCashTCH'
$ws.Range("B19").Value = 'BitKan'
$ws.Range("C19").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.001055'
$ws.Range("E19").Value = '18BitKanKAN'
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.004115'
$ws.Range("E20").Value = '19HotbitTokenHTB'
$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0001502'
$ws.Range("E21").Value = '20NitroExNTX'
$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.967'
$ws.Range("E22").Value = '21LEOLEO'
$ws.Range("B23").Value = 'GateToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.379'
$ws.Range("E23").Value = '22GateTokenGT'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.100'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3272'
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'
$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1312'
$ws.Range("E26").Value = '25ProBitTokenPROBBestin24h'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04184'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006869'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003504'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009827'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005636'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6809'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.02472'
$ws.Range("E48").Value = '47BOLOBOLOWorstin24h'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002103'
